$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.216178999999999
$ws.Range("H2").Value = 12.648537
$ws.Range("I2").Value = 0.01683192247764961
$ws.Range("J2").Value = 0.01683192247764961
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 88.88470006331298
$ws.Range("R2").Value = 799.962300569817
$ws.Range("S2").Value = 0.0009621228125480959
$ws.Range("T2").Value = 0.0009621228125480961

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.216178999999999
$ws.Range("H3").Value = 12.648537
$ws.Range("I3").Value = 0.01683192247764961
$ws.Range("J3").Value = 0.01683192247764961
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 1271.600213205278
$ws.Range("R3").Value = 11444.4019188475
$ws.Range("S3").Value = 0.01376429883539418
$ws.Range("T3").Value = 0.01376429883539419

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.216178999999999
$ws.Range("H4").Value = 12.648537
$ws.Range("I4").Value = 0.01683192247764961
$ws.Range("J4").Value = 0.01683192247764961
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 194.514470804357
$ws.Range("R4").Value = 1750.630237239213
$ws.Range("S4").Value = 0.002105500829707328
$ws.Range("T4").Value = 0.002105500829707328

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 212.7693433333334
$ws.Range("H5").Value = 638.30803
$ws.Range("I5").Value = 0.8494224492382987
$ws.Range("J5").Value = 0.8494224492382987
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 4485.563650132359
$ws.Range("R5").Value = 40370.07285119123
$ws.Range("S5").Value = 0.04855349809196388
$ws.Range("T5").Value = 0.04855349809196388

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 212.7693433333334
$ws.Range("H6").Value = 638.30803
$ws.Range("I6").Value = 0.8494224492382987
$ws.Range("J6").Value = 0.8494224492382987
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 64171.26558104238
$ws.Range("R6").Value = 577541.3902293814
$ws.Range("S6").Value = 0.6946149166462302
$ws.Range("T6").Value = 0.6946149166462302

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 212.7693433333334
$ws.Range("H7").Value = 638.30803
$ws.Range("I7").Value = 0.8494224492382987
$ws.Range("J7").Value = 0.8494224492382987
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 9816.166776096055
$ws.Range("R7").Value = 88345.50098486448
$ws.Range("S7").Value = 0.1062540345001046
$ws.Range("T7").Value = 0.1062540345001046

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 33.501551
$ws.Range("H8").Value = 100.504653
$ws.Range("I8").Value = 0.1337456282840517
$ws.Range("J8").Value = 0.1337456282840517
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 706.2734557263303
$ws.Range("R8").Value = 6356.461101536972
$ws.Range("S8").Value = 0.007644980555342522
$ws.Range("T8").Value = 0.007644980555342522

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 33.501551
$ws.Range("H9").Value = 100.504653
$ws.Range("I9").Value = 0.1337456282840517
$ws.Range("J9").Value = 0.1337456282840517
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 10104.07276216391
$ws.Range("R9").Value = 90936.65485947522
$ws.Range("S9").Value = 0.1093704416755548
$ws.Range("T9").Value = 0.1093704416755548

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 33.501551
$ws.Range("H10").Value = 100.504653
$ws.Range("I10").Value = 0.1337456282840517
$ws.Range("J10").Value = 0.1337456282840517
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 1545.602419605567
$ws.Range("R10").Value = 13910.4217764501
$ws.Range("S10").Value = 0.01673020605315437
$ws.Range("T10").Value = 0.01673020605315437
